# Arreglo el crud y tablas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three extra sample rows (old rows 4,5,6 - "jesucita", "sadfasd", the stray row)
$ws.Rows("4:6").Delete()

# Remove the "nombre" and "edad" columns (old B and C); this shifts
# dia/fecha/mes/unidad_numero/clase_numero/curso/caracter_clase/
# contenidos_tematicos/actividades/observaciones left into B..K,
# which already matches the desired header order.
$ws.Range("B1:C1").EntireColumn.Delete()

# Fill in the two remaining data rows (id=1 and id=2) with the new
# dia/fecha/mes/unidad_numero/clase_numero/curso/caracter_clase/
# contenidos_tematicos/actividades/observaciones sample values.
# Force the cells to Text first so numeric-looking values ("23", "2",
# "09"...) and the date-looking value keep their original literal text
# instead of being auto-converted to a number/date by Excel.
$ws.Range("B2:K3").NumberFormat = "@"

$row2 = New-Object 'object[,]' 1,10
$row2[0,0] = "23"
$row2[0,1] = "2024-09-07"
$row2[0,2] = "2"
$row2[0,3] = "5"
$row2[0,4] = "6"
$row2[0,5] = "5b"
$row2[0,6] = "aaaaa"
$row2[0,7] = "nnnnnn"
$row2[0,8] = "sssss"
$row2[0,9] = "ccccc"
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object 'object[,]' 1,10
$row3[0,0] = "5"
$row3[0,1] = "2024-09-07"
$row3[0,2] = "23"
$row3[0,3] = "3"
$row3[0,4] = "1"
$row3[0,5] = "5b"
$row3[0,6] = "asdfa"
$row3[0,7] = "sdf"
$row3[0,8] = "sda"
$row3[0,9] = "fda"
$ws.Range("B3:K3").Value = $row3

# Drop the temporary "@" number format again now that the values are
# locked in as text, so the cells end up back on the default (unstyled)
# look, matching the original plain data rows.
$ws.Range("B2:K3").ClearFormats()
